# Actualización automática 2025-08-05 16:45:08
#
# Updates the per-advisor sales figures for HIDALGO HIDALGO PEDRO GUSTAVO
# across the three sheets of the workbook:
#   - "VENTAS POR GRUPO"      : sales by product group (row 7 = client row,
#                                row 22 = totals row)
#   - "VENTA MENSUAL"         : monthly sales (row 7 = client row,
#                                row 22 = totals row, column F = agosto)
#   - "CUMPLIMIENTO MENSUAL"  : budget vs sales compliance (rows 7/8 = product
#                                rows, row 19 = TOTAL row)

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H7").Value = 68.73
$ws1.Range("I7").Value = 42.63
$ws1.Range("H22").Value = "1 de 20"
$ws1.Range("I22").Value = "1 de 20"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 111.36
$ws2.Range("F22").Value = 111.36

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# INODOROS
$ws3.Range("D7").Value = 68.73
$ws3.Range("E7").Value = 2331.27
$ws3.Range("F7").Value = 0.0286375

# LAVABOS
$ws3.Range("D8").Value = 42.63
$ws3.Range("E8").Value = 582.37
$ws3.Range("F8").Value = 0.068208

# TOTAL
$ws3.Range("D19").Value = 111.36
$ws3.Range("E19").Value = 59276.86762291769
$ws3.Range("F19").Value = 0.001875119101163858
